$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
# Remove the 2024-09-21 event ("南宁·小野爷爷&娃展2.0"); every subsequent row
# shifts up by one and the refreshed scrape bumped a few "want to go" counts.
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Rows.Item(2).Delete()

$ws1.Range("A2").Value = 1
$ws1.Range("F2").Value = 408

$ws1.Range("A3").Value = 2
$ws1.Range("F3").Value = 5114

$ws1.Range("A4").Value = 3
$ws1.Range("F4").Value = 44

$ws1.Range("A5").Value = 4
$ws1.Range("F5").Value = 46

$ws1.Range("A6").Value = 5
$ws1.Range("F6").Value = 47

$ws1.Range("A7").Value = 6
$ws1.Range("F7").Value = 504

# --- Sheet "全部类型" (all types) ---
# Same underlying change applied to the combined listing sheet.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Rows.Item(2).Delete()

$ws4.Range("A2").Value = 1
$ws4.Range("F2").Value = 408

$ws4.Range("A3").Value = 2
$ws4.Range("F3").Value = 5114

$ws4.Range("A4").Value = 3
$ws4.Range("F4").Value = 6

$ws4.Range("A5").Value = 4
$ws4.Range("F5").Value = 44

$ws4.Range("A6").Value = 5
$ws4.Range("F6").Value = 46

$ws4.Range("A7").Value = 6
$ws4.Range("F7").Value = 7

$ws4.Range("A8").Value = 7
$ws4.Range("F8").Value = 47

$ws4.Range("A9").Value = 8
$ws4.Range("F9").Value = 504
